$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cluster / cell-count data for rows 2-23 (row 1 is the header and is unchanged)
$data = @(
    @("MeV.Fib.0", 456),
    @("MeV.VLMC.0", 345),
    @("MeV.Endothelial.0", 279),
    @("MeV.FibCollagen.3", 255),
    @("MeV.FibUnknown.8", 250),
    @("MeV.Endothelial.2", 246),
    @("MeV.Fib.2", 229),
    @("MeV.Fib.1", 205),
    @("MeV.FibCollagen.0", 198),
    @("MeV.Pericytes.0", 189),
    @("MeV.Endothelial.1", 115),
    @("MeV.VLMC.1", 112),
    @("MeV.FibLaminin.0", 90),
    @("MeV.Endothelial.3", 88),
    @("MeV.Fib.4", 84),
    @("MeV.Epithelial.0", 72),
    @("MeV.FibCollagen.2", 54),
    @("MeV.SMC.0", 49),
    @("MeV.FibCollagen.1", 47),
    @("MeV.Fib.3", 40),
    @("MeV.EndoUnknowed.4", 32),
    @("MeV.FibProlif.0", 21)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# The previous sheet had 24 data rows (A1:B24); the updated sheet only has 23
# (A1:B23), so the last row that is no longer present must be removed.
$ws.Cells.Item(24, 1).EntireRow.Delete()
